$d = $word.ActiveDocument

# The document body is a pair of near-identical "spec sheet" blocks, each a
# run of text separated by manual line breaks (<w:br/>) inside a single
# paragraph. Right after each block's "Utilities Req'd: ..." line there is
# a blank line made of two consecutive manual line breaks. This change
# turns that blank line into an actual paragraph break: the first <w:br/>
# run is removed and the paragraph is split there, so the second <w:br/>
# becomes the leading run of the new paragraph.
#
# (The very end of the document also happens to contain extra manual line
# breaks in a row, but those are unrelated boilerplate and must be left
# alone, so we anchor the search on "Utilities Req" and only look for the
# double line-break that immediately follows it.)

$lineBreak = [string][char]11
$pattern = $lineBreak + $lineBreak
$label = "Utilities Req"

$labelSearchFrom = -1
while ($true) {
    $t = $d.Content.Text
    $labelIdx = $t.IndexOf($label, $labelSearchFrom + 1)
    if ($labelIdx -lt 0) {
        break
    }
    $labelSearchFrom = $labelIdx

    $idx = $t.IndexOf($pattern, $labelIdx)
    if ($idx -lt 0) {
        continue
    }

    # Split the paragraph between the two manual line breaks: the range
    # covers just the second break character, so InsertParagraphBefore
    # places the new paragraph mark exactly between break #1 and break #2.
    $splitRange = $d.Range($idx + 1, $idx + 2)
    $splitRange.InsertParagraphBefore()

    # The first manual line break (now redundant, since the paragraph mark
    # performs that job) is deleted outright.
    $firstBreakRange = $d.Range($idx, $idx + 1)
    $firstBreakRange.Text = ""
}
